$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.083.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "'2.483.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'585.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "'171.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.66%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.513"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").Value = "'2.481.55"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").Value = "'4.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").Value = "'2.923.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").Value = "'25.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("D16").Value = "'66.972.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").Value = "'2.498.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.19%  "
$ws.Range("D19").Value = "'10.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.90%  "
$ws.Range("E20").Value = "  -4.45%  "
$ws.Range("D21").Value = "'350.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.52%  "
$ws.Range("D22").Value = "'4.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'68.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.67%  "
$ws.Range("D25").Value = "'4.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.57%  "
$ws.Range("D26").Value = "'1.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.36%  "
$ws.Range("D27").Value = "'9.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("D28").Value = "'0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.85%  "
$ws.Range("D30").Value = "'0.0₃0900"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("D31").Value = "'507.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("D32").Value = "'7.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.40%  "
$ws.Range("E33").Value = "  -3.64%  "
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "'159.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.77%  "
$ws.Range("E37").Value = "  -7.52%  "
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("D39").Value = "'18.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.94%  "
$ws.Range("E40").Value = "  -6.03%  "
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").Value = "'1.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.59%  "
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("D44").Value = "'4.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.76%  "
$ws.Range("D45").Value = "'2.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("D46").Value = "'38.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").Value = "'142.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("E48").Value = "  -4.13%  "
$ws.Range("D49").Value = "'3.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.46%  "
$ws.Range("D50").Value = "'0.0₆0250"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.54%  "
$ws.Range("D51").Value = "'0.0728"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.32%  "
